# Insert a new "account" worksheet between "assumption" and "result",
# populate it with the account register, and leave it as the active sheet
# (matching the commit "update excel model with accounts").

$wb = $excel.ActiveWorkbook

$assumption = $wb.Worksheets.Item("assumption")
$ws = $wb.Worksheets.Add($null, $assumption)
$ws.Name = "account"

# Header row
$ws.Range("A1").Value = "名称"
$ws.Range("B1").Value = "类型"
$ws.Range("C1").Value = "余额"

# Type column filled first (matches original authoring order)
$ws.Range("B2").Value = "本金帐"
$ws.Range("B3").Value = "利息帐"

# Name column
$ws.Range("A2").Value = "账户P"
$ws.Range("A3").Value = "账户I"

# Balances
$ws.Range("C2").Value = 1200
$ws.Range("C3").Value = 300

# Leave the selection where the author left it
[void]$ws.Range("C4").Select()
